# Update countries & provincias Spain
# - Reorders "Pakistan" ahead of "Japon" (and refreshes Pakistan's figures)
# - Reorders "Laos" ahead of "Santa Lucia" (and refreshes Laos' figures),
#   shifting Santa Lucia / Sudan / Liberia / Curazao down one row
# - Refreshes case counts for several other countries
# - Bumps the "datos actualizados" timestamp from 09:52 to 10:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 10:22"

# --- Row 13: Suiza ------------------------------------------------------
$ws.Range("B13").Value = 21712
$ws.Range("C13").Value = 55
$ws.Range("E13").Value = 12887
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 769

# --- Row 17: Austria ------------------------------------------------------
$ws.Range("B17").Value = 12390
$ws.Range("C17").Value = 93
$ws.Range("E17").Value = 8707

# --- Row 21: Israel ------------------------------------------------------
$ws.Range("D21").Value = 683
$ws.Range("E21").Value = 8264

# --- Row 24: Australia ------------------------------------------------------
$ws.Range("B24").Value = 5908
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 2547
$ws.Range("E24").Value = 3315
$ws.Range("F24").Value = 93

# --- Row 25: Noruega ------------------------------------------------------
$ws.Range("B25").Value = 5866
$ws.Range("C25").Value = 1
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 78

# --- Row 31: Polonia ------------------------------------------------------
$ws.Range("B31").Value = 4532
$ws.Range("C31").Value = 119
$ws.Range("E31").Value = 4230
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 111

# --- Rows 33-34: Pakistan now listed before Japon ------------------------
$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 3918
$ws.Range("C33").Value = 152
$ws.Range("D33").Value = 429
$ws.Range("E33").Value = 3435
$ws.Range("F33").Value = 28
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 54

$ws.Range("A34").Value = "Japon"
$ws.Range("B34").Value = 3906
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 592
$ws.Range("E34").Value = 3222
$ws.Range("F34").Value = 79
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 92

# --- Row 71: Bosnia y Herzegovina ------------------------------------------------------
$ws.Range("B71").Value = 740
$ws.Range("C71").Value = 66
$ws.Range("E71").Value = 642

# --- Row 75: Kuwait ------------------------------------------------------
$ws.Range("D75").Value = 105
$ws.Range("E75").Value = 559

# --- Rows 174-178: Laos now listed before Santa Lucia ----------------------
# (Santa Lucia, Sudan, Liberia, Curazao each shift down one row)
$ws.Range("A174").Value = "Laos"
$ws.Range("C174").Value = 2
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 14

$ws.Range("A175").Value = "Santa Lucia"
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 1
$ws.Range("E175").Value = 13
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = "Sudan"
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 2
$ws.Range("E176").Value = 10
$ws.Range("H176").Value = 2

$ws.Range("A177").Value = "Liberia"
$ws.Range("B177").Value = 14
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 8
$ws.Range("H177").Value = 3

$ws.Range("A178").Value = "Curazao"
$ws.Range("B178").Value = 13
$ws.Range("D178").Value = 5
$ws.Range("E178").Value = 7
$ws.Range("H178").Value = 1
